$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Goal (per the diff): split the existing single paragraph
#     "New line" + <bookmarkStart _GoBack/><bookmarkEnd/>
# into two paragraphs:
#     "New line"
#     "New line 2" + <bookmarkStart _GoBack/><bookmarkEnd/>
# i.e. insert a new paragraph break + "New line 2" right after "New line",
# and the (hidden) "_GoBack" bookmark should end up collapsed at the end of
# the new "New line 2" paragraph (where it originally sat relative to the
# end of the text), instead of staying in the first paragraph.
# ---------------------------------------------------------------------------

# 1) Find the existing "New line" text and collapse the range to its end -
#    this is where the new paragraph break must be inserted.
$target = $d.Content
$null = $target.Find.Execute("New line", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Collapse(0)  # wdCollapseEnd

# 2) Drop the current "_GoBack" bookmark (it is hidden from the Bookmarks
#    collection/Count, but still addressable by name). It will be recreated
#    after the new paragraph/text have been inserted, in the right spot.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3) Insert the paragraph break and the new line's text right after
#    "New line". A temporary trailing marker character ("#") is appended
#    after "New line 2" and removed afterwards - this works around a quirk
#    where adding a bookmark via a collapsed range sitting exactly at
#    "end of paragraph text" (immediately before the paragraph mark) gets
#    mis-resolved; having extra content after the insertion point keeps the
#    anchor unambiguous while we (re)create the bookmark.
$target.InsertAfter([char]13 + "New line 2" + "#")

# 4) Re-add the "_GoBack" bookmark, collapsed, right after "New line 2" -
#    i.e. immediately before the temporary marker character.
$markerPos = $target.End - 1
$bmRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 5) Remove the temporary marker character.
$d.Range($markerPos, $markerPos + 1).Delete()
